# Update "想去人数" (want-to-go count) figures to match the regenerated
# gh-pages output at commit 456a3b4.
#
# 展览 sheet: row 3 (环形宇宙动漫游戏嘉年华) F3 2781 -> 2786
#             row 4 (MAX特摄同人only2.0)     F4 126  -> 127
# 全部类型 sheet: same two events appear again at rows 7/8
#             F7 2781 -> 2786
#             F8 126  -> 127

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 2786
$wsExpo.Range("F4").Value = 127

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2786
$wsAll.Range("F8").Value = 127
